# Weekly update: prepend a new week's Limón price record for
# "Feria Lagunitas de Puerto Montt" ahead of the existing history,
# pushing the previous rows down by two and growing the used range
# from A1:T499 to A1:T501.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two fresh rows right above the first data row of this block
# (row 469), shifting every existing row (469-499) down to (471-501).
$ws.Range("A469:A470").EntireRow.Insert()

# New row 469: "1a plateado"
$ws.Cells.Item(469, 1).Value  = 4
$ws.Cells.Item(469, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(469, 3).Value  = "Los Lagos"
$ws.Cells.Item(469, 4).Value  = 44714
$ws.Cells.Item(469, 5).Value  = 10
$ws.Cells.Item(469, 6).Value  = "Fruta"
$ws.Cells.Item(469, 7).Value  = 100102
$ws.Cells.Item(469, 8).Value  = "Cítricos"
$ws.Cells.Item(469, 9).Value  = 100102003
$ws.Cells.Item(469, 10).Value = "Limón"
$ws.Cells.Item(469, 11).Value = "Sin especificar"
$ws.Cells.Item(469, 12).Value = "1a plateado"
$ws.Cells.Item(469, 13).Value = 800
$ws.Cells.Item(469, 14).Value = 12000
$ws.Cells.Item(469, 15).Value = 13000
$ws.Cells.Item(469, 16).Value = 12500
$ws.Cells.Item(469, 17).Value = "$/malla 18 kilos"
$ws.Cells.Item(469, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(469, 19).Value = 694
$ws.Cells.Item(469, 20).Value = 18

# New row 470: "2a plateado"
$ws.Cells.Item(470, 1).Value  = 4
$ws.Cells.Item(470, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(470, 3).Value  = "Los Lagos"
$ws.Cells.Item(470, 4).Value  = 44714
$ws.Cells.Item(470, 5).Value  = 10
$ws.Cells.Item(470, 6).Value  = "Fruta"
$ws.Cells.Item(470, 7).Value  = 100102
$ws.Cells.Item(470, 8).Value  = "Cítricos"
$ws.Cells.Item(470, 9).Value  = 100102003
$ws.Cells.Item(470, 10).Value = "Limón"
$ws.Cells.Item(470, 11).Value = "Sin especificar"
$ws.Cells.Item(470, 12).Value = "2a plateado"
$ws.Cells.Item(470, 13).Value = 300
$ws.Cells.Item(470, 14).Value = 10000
$ws.Cells.Item(470, 15).Value = 10000
$ws.Cells.Item(470, 16).Value = 10000
$ws.Cells.Item(470, 17).Value = "$/malla 18 kilos"
$ws.Cells.Item(470, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(470, 19).Value = 556
$ws.Cells.Item(470, 20).Value = 18
